$wb = $excel.ActiveWorkbook

# --- TravelInsurance sheet (sheet1) ---
$ws1 = $wb.Worksheets.Item("TravelInsurance")

# Fix header typo: "phonenumber" -> "phoneNumber"
$ws1.Range("E1").Value = "phoneNumber"

# Apply explicit number formats to the data row so ages are numeric and the
# phone number is stored/displayed as text (keeps leading formatting stable)
$ws1.Range("C2:D2").NumberFormat = "0"
$ws1.Range("E2").NumberFormat = "@"

# Set up printing for this sheet
$ws1.PageSetup.Orientation = 1

# Make TravelInsurance the active/selected sheet with F9 as the active cell
$ws1.Activate()
$ws1.Range("F9").Select()
